$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2 - Numbers")

# Add new column AA with values 100-129 across rows 1-30 (matches reader/sheet2)
for ($i = 0; $i -lt 30; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 27).Value2 = 100 + $i
}

# Make Sheet2 the active sheet/tab with AA1:AA30 selected (AA1 active cell)
$ws2.Activate()
$ws2.Range("AA1:AA30").Select()

# Try to scroll the view so column O is the leftmost visible column
$excel.ActiveWindow.ScrollColumn = 15
$excel.ActiveWindow.ScrollRow = 1

# Sheet4: pageSetup paperSize 0 -> 9 (Letter)
$ws4 = $wb.Worksheets.Item("Sheet4 - Dates")
$ws4.PageSetup.PaperSize = 9
